$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")
$ws.Activate()

# The only real data edit: RFduvar gücü (kW) [Pdu] (C7) changed from 36 to 10.
$ws.Range("C7").Value = 10

# View/selection state: scroll so row 4 is at the top and select C8
# (matches topLeftCell="A4" / selection activeCell="C8" in the target file).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select() | Out-Null
